$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.06449866666666666
$ws.Range("N2").Value = 0.193496
$ws.Range("O2").Value = 0.004525829983623641
$ws.Range("P2").Value = 0.004525829983623642
$ws.Range("Q2").Value = 0.5906116677208889
$ws.Range("R2").Value = 5.315505009488001
$ws.Range("S2").Value = 0.0043872114682169
$ws.Range("T2").Value = 0.004387211468216901

# Row 3
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("O3").Value = 0.745188142173877
$ws.Range("P3").Value = 0.7451881421738772
$ws.Range("S3").Value = 0.7223642901200832
$ws.Range("T3").Value = 0.7223642901200834

# Row 4
$ws.Range("I4").Value = 0.9693716918425304
$ws.Range("J4").Value = 0.9693716918425304
$ws.Range("M4").Value = 3.566885000000001
$ws.Range("N4").Value = 10.700655
$ws.Range("O4").Value = 0.2502860278424993
$ws.Range("P4").Value = 0.2502860278424993
$ws.Range("Q4").Value = 32.66182089167668
$ws.Range("R4").Value = 293.9563880250901
$ws.Range("S4").Value = 0.2426201902542302
$ws.Range("T4").Value = 0.2426201902542302

# Row 5
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.06449866666666666
$ws.Range("N5").Value = 0.193496
$ws.Range("O5").Value = 0.004525829983623641
$ws.Range("P5").Value = 0.004525829983623642
$ws.Range("Q5").Value = 0.01866099073511111
$ws.Range("R5").Value = 0.167948916616
$ws.Range("S5").Value = 0.0001386185154067406
$ws.Range("T5").Value = 0.0001386185154067406

# Row 6
$ws.Range("G6").Value = 0.2893236666666667
$ws.Range("H6").Value = 0.867971
$ws.Range("I6").Value = 0.03062830815746963
$ws.Range("J6").Value = 0.03062830815746962
$ws.Range("O6").Value = 0.745188142173877
$ws.Range("P6").Value = 0.7451881421738772
$ws.Range("Q6").Value = 3.072574327214889
$ws.Range("R6").Value = 27.653168944934
$ws.Range("S6").Value = 0.02282385205379379
$ws.Range("T6").Value = 0.0228238520537938

# Row 7
$ws.Range("G7").Value = 0.2893236666666667
$ws.Range("H7").Value = 0.867971
$ws.Range("I7").Value = 0.03062830815746963
$ws.Range("J7").Value = 0.03062830815746962
$ws.Range("M7").Value = 3.566885000000001
$ws.Range("N7").Value = 10.700655
$ws.Range("O7").Value = 0.2502860278424993
$ws.Range("P7").Value = 0.2502860278424993
$ws.Range("Q7").Value = 1.031984246778334
$ws.Range("R7").Value = 9.287858221005001
$ws.Range("S7").Value = 0.007665837588269091
$ws.Range("T7").Value = 0.007665837588269091

